# Insert a new data row at row 541 (pushes existing rows 541:656 down to 542:657)
# and populate it with the new "Perejil" observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("541:541").Insert()

$ws.Range("A541").Value2 = 9
$ws.Range("B541").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C541").Value2 = "Metropolitana"
$ws.Range("D541").Value2 = 45275
$ws.Range("E541").Value2 = 13
$ws.Range("F541").Value2 = 100112044
$ws.Range("G541").Value2 = "Perejil"
$ws.Range("H541").Value2 = "Sin especificar"
$ws.Range("I541").Value2 = "Primera"
$ws.Range("J541").Value2 = 70
$ws.Range("K541").Value2 = 10000
$ws.Range("L541").Value2 = 12000
$ws.Range("M541").Value2 = 11000
$ws.Range("N541").Value2 = "`$/docena de atados"
$ws.Range("O541").Value2 = "Región Metropolitana"
$ws.Range("P541").Value2 = 3667
$ws.Range("Q541").Value2 = 3
$ws.Range("R541").Value2 = "Hortaliza"
